$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab query (B2) to remove the trailing Cohort line.
# This is the only unique-string cell that changes; the shared-string table
# compaction + reindexing that Excel performs on save automatically updates
# every other cell that referenced a shifted shared-string index (C2:E2, B3:E3, B4:E4).
$casesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN [''UBC01''] and demo.breed in [''Mixed Breed'', ''Scottish Terrier'',''Shetland Sheepdog'']and diag.disease_term in [''Bladder Cancer'',''Healthy Control''] and diag.primary_disease_site in [''Bladder'', ''Bladder, Prostate'', ''Bladder, Urethra'', ''Bladder, Urethra, Prostate'']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'
$ws.Range("B2").Value = $casesQuery

# Row heights: row 2 shrinks now that the Cohort line is gone.
$ws.Rows.Item(2).RowHeight = 304.5

# Selection / scroll position moved from D4 to B2.
$ws.Range("B2").Select() | Out-Null

